# Updated activity till excel form
# Re-sort/shuffle the per-match stat rows (runs, balls, fours, sixes) for
# Nicholas Pooran while keeping playerName/teamName columns untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns C (runs), D (balls), E (fours), F (sixes)
$data = @{
    2  = @("22", "10", "0", "3")
    3  = @("24", "12", "2", "2")
    4  = @("2",  "6",  "0", "0")
    5  = @("32", "28", "2", "0")
    6  = @("2",  "3",  "0", "0")
    7  = @("16", "10", "2", "1")
    8  = @("17", "18", "1", "0")
    9  = @("53", "28", "6", "3")
    10 = @("6",  "1",  "0", "1")
    11 = @("25", "8",  "1", "3")
    12 = @("44", "27", "3", "2")
    13 = @("0",  "3",  "0", "0")
    14 = @("77", "37", "5", "7")
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("C$row").Value = $values[0]
    $ws.Range("D$row").Value = $values[1]
    $ws.Range("E$row").Value = $values[2]
    $ws.Range("F$row").Value = $values[3]
}
